# Cards changes. Hound, Ogre, Tomb, Shopkeeper.
# Insert the new "Tomb" (墓碑) trap card as row 2, pushing the existing
# rows (Quicksand..Remote bomb) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 2 (shifts rows 2:12 -> 3:13)
$ws.Rows(2).Insert()

# --- Populate the new row 2 with the Tomb card data ---
$ws.Cells.Item(2, 1).Value2 = "墓碑"
$ws.Cells.Item(2, 2).Value2 = 1
$ws.Cells.Item(2, 3).Value2 = 3
$ws.Cells.Item(2, 4).Value2 = "回合结束时：将主牌堆顶2张牌送墓。<br>`n开战时：用墓地顶端1张牌替换本牌。"
$ws.Cells.Item(2, 4).WrapText = $true
$ws.Cells.Item(2, 5).Value2 = "Tomb"

# Row height for the new wrapped 2-line effect text
$ws.Rows(2).RowHeight = 28.5

# --- Restore the view: selection moves to D6 (also resets scroll to top) ---
$ws.Range("D6").Select()
